$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.170692920684814
$ws.Range("B1").Value = 2.381422996520996
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.377081632614136
$ws.Range("E1").Value = 1.210865020751953
